$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value2 = "ECs"
$ws.Cells.Item(2, 2).Value2 = "Efnb3"
$ws.Cells.Item(2, 3).Value2 = "Rhbdl2"
$ws.Cells.Item(2, 4).Value2 = "ECs"
$ws.Cells.Item(2, 5).Value2 = 1
$ws.Cells.Item(2, 6).Value2 = 0.3333333333333333
$ws.Cells.Item(2, 7).Value2 = 0.071327
$ws.Cells.Item(2, 8).Value2 = 0.213981
$ws.Cells.Item(2, 9).Value2 = 0.1064107741026886
$ws.Cells.Item(2, 10).Value2 = 0.1213845878939322
$ws.Cells.Item(2, 11).Value2 = 3
$ws.Cells.Item(2, 12).Value2 = 1
$ws.Cells.Item(2, 13).Value2 = 1.806704333333333
$ws.Cells.Item(2, 14).Value2 = 5.420113000000001
$ws.Cells.Item(2, 15).Value2 = 0.9475373327607437
$ws.Cells.Item(2, 16).Value2 = 0.9644023730209256
$ws.Cells.Item(2, 17).Value2 = 0.1288667999836667
$ws.Cells.Item(2, 18).Value2 = 1.159801199853
$ws.Cells.Item(2, 19).Value2 = 0.1008281810702675
$ws.Cells.Item(2, 20).Value2 = 0.1170635846130754

# Row 3
$ws.Cells.Item(3, 1).Value2 = "ECs"
$ws.Cells.Item(3, 2).Value2 = "Efnb3"
$ws.Cells.Item(3, 3).Value2 = "Rhbdl2"
$ws.Cells.Item(3, 4).Value2 = "MuSCs"
$ws.Cells.Item(3, 5).Value2 = 1
$ws.Cells.Item(3, 6).Value2 = 0.3333333333333333
$ws.Cells.Item(3, 7).Value2 = 0.071327
$ws.Cells.Item(3, 8).Value2 = 0.213981
$ws.Cells.Item(3, 9).Value2 = 0.1064107741026886
$ws.Cells.Item(3, 10).Value2 = 0.1213845878939322
$ws.Cells.Item(3, 11).Value2 = 1
$ws.Cells.Item(3, 12).Value2 = 0.5
$ws.Cells.Item(3, 13).Value2 = 0.1000325
$ws.Cells.Item(3, 14).Value2 = 0.200065
$ws.Cells.Item(3, 15).Value2 = 0.05246266723925631
$ws.Cells.Item(3, 16).Value2 = 0.03559762697907432
$ws.Cells.Item(3, 17).Value2 = 0.0071350181275
$ws.Cells.Item(3, 18).Value2 = 0.042810108765
$ws.Cells.Item(3, 19).Value2 = 0.005582593032421024
$ws.Cells.Item(3, 20).Value2 = 0.00432100328085686

# Row 4
$ws.Cells.Item(4, 1).Value2 = "FAPs"
$ws.Cells.Item(4, 2).Value2 = "Efnb3"
$ws.Cells.Item(4, 3).Value2 = "Rhbdl2"
$ws.Cells.Item(4, 4).Value2 = "ECs"
$ws.Cells.Item(4, 5).Value2 = 2
$ws.Cells.Item(4, 6).Value2 = 0.6666666666666666
$ws.Cells.Item(4, 7).Value2 = 0.3509106666666666
$ws.Cells.Item(4, 8).Value2 = 1.052732
$ws.Cells.Item(4, 9).Value2 = 0.5235138962929958
$ws.Cells.Item(4, 10).Value2 = 0.5971812449832231
$ws.Cells.Item(4, 11).Value2 = 3
$ws.Cells.Item(4, 12).Value2 = 1
$ws.Cells.Item(4, 13).Value2 = 1.806704333333333
$ws.Cells.Item(4, 14).Value2 = 5.420113000000001
$ws.Cells.Item(4, 15).Value2 = 0.9475373327607437
$ws.Cells.Item(4, 16).Value2 = 0.9644023730209256
$ws.Cells.Item(4, 17).Value2 = 0.6339918220795556
$ws.Cells.Item(4, 18).Value2 = 5.705926398716
$ws.Cells.Item(4, 19).Value2 = 0.4960489609566499
$ws.Cells.Item(4, 20).Value2 = 0.5759230097854111

# Row 5
$ws.Cells.Item(5, 1).Value2 = "FAPs"
$ws.Cells.Item(5, 2).Value2 = "Efnb3"
$ws.Cells.Item(5, 3).Value2 = "Rhbdl2"
$ws.Cells.Item(5, 4).Value2 = "MuSCs"
$ws.Cells.Item(5, 5).Value2 = 2
$ws.Cells.Item(5, 6).Value2 = 0.6666666666666666
$ws.Cells.Item(5, 7).Value2 = 0.3509106666666666
$ws.Cells.Item(5, 8).Value2 = 1.052732
$ws.Cells.Item(5, 9).Value2 = 0.5235138962929958
$ws.Cells.Item(5, 10).Value2 = 0.5971812449832231
$ws.Cells.Item(5, 11).Value2 = 1
$ws.Cells.Item(5, 12).Value2 = 0.5
$ws.Cells.Item(5, 13).Value2 = 0.1000325
$ws.Cells.Item(5, 14).Value2 = 0.200065
$ws.Cells.Item(5, 15).Value2 = 0.05246266723925631
$ws.Cells.Item(5, 16).Value2 = 0.03559762697907432
$ws.Cells.Item(5, 17).Value2 = 0.03510247126333333
$ws.Cells.Item(5, 18).Value2 = 0.21061482758
$ws.Cells.Item(5, 19).Value2 = 0.02746493533634598
$ws.Cells.Item(5, 20).Value2 = 0.02125823519781198

# Row 6
$ws.Cells.Item(6, 1).Value2 = "MuSCs"
$ws.Cells.Item(6, 2).Value2 = "Efnb3"
$ws.Cells.Item(6, 3).Value2 = "Rhbdl2"
$ws.Cells.Item(6, 4).Value2 = "ECs"
$ws.Cells.Item(6, 5).Value2 = 1
$ws.Cells.Item(6, 6).Value2 = 0.5
$ws.Cells.Item(6, 7).Value2 = 0.248061
$ws.Cells.Item(6, 8).Value2 = 0.496122
$ws.Cells.Item(6, 9).Value2 = 0.3700753296043157
$ws.Cells.Item(6, 10).Value2 = 0.2814341671228447
$ws.Cells.Item(6, 11).Value2 = 3
$ws.Cells.Item(6, 12).Value2 = 1
$ws.Cells.Item(6, 13).Value2 = 1.806704333333333
$ws.Cells.Item(6, 14).Value2 = 5.420113000000001
$ws.Cells.Item(6, 15).Value2 = 0.9475373327607437
$ws.Cells.Item(6, 16).Value2 = 0.9644023730209256
$ws.Cells.Item(6, 17).Value2 = 0.4481728836310001
$ws.Cells.Item(6, 18).Value2 = 2.689037301786001
$ws.Cells.Item(6, 19).Value2 = 0.3506601907338264
$ws.Cells.Item(6, 20).Value2 = 0.2714157786224392

# Row 7
$ws.Cells.Item(7, 1).Value2 = "MuSCs"
$ws.Cells.Item(7, 2).Value2 = "Efnb3"
$ws.Cells.Item(7, 3).Value2 = "Rhbdl2"
$ws.Cells.Item(7, 4).Value2 = "MuSCs"
$ws.Cells.Item(7, 5).Value2 = 1
$ws.Cells.Item(7, 6).Value2 = 0.5
$ws.Cells.Item(7, 7).Value2 = 0.248061
$ws.Cells.Item(7, 8).Value2 = 0.496122
$ws.Cells.Item(7, 9).Value2 = 0.3700753296043157
$ws.Cells.Item(7, 10).Value2 = 0.2814341671228447
$ws.Cells.Item(7, 11).Value2 = 1
$ws.Cells.Item(7, 12).Value2 = 0.5
$ws.Cells.Item(7, 13).Value2 = 0.1000325
$ws.Cells.Item(7, 14).Value2 = 0.200065
$ws.Cells.Item(7, 15).Value2 = 0.05246266723925631
$ws.Cells.Item(7, 16).Value2 = 0.03559762697907432
$ws.Cells.Item(7, 17).Value2 = 0.0248141619825
$ws.Cells.Item(7, 18).Value2 = 0.09925664793
$ws.Cells.Item(7, 19).Value2 = 0.01941513887048932
$ws.Cells.Item(7, 20).Value2 = 0.01001838850040549

